# Refresh the cryptos price list: GitHub Actions re-scraped coinranking.com
# and the sheet's Price (D) / Volume(1h) (E) columns need the new readings.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.911.46"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "2.305.32"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'305.97"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").Value = "'97.17"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").Value = "'0.512"
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.503"
$ws.Range("E9").Value = "  -2.66%  "
$ws.Range("D10").Value = "'35.60"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "'0.0795"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "'18.27"
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "'6.76"
$ws.Range("E14").Value = "  -1.82%  "
$ws.Range("D15").Value = "2.663.87"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("D16").Value = "2.312.85"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("D17").Value = "'0.783"
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "42.852.72"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").Value = "'13.09"
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").Value = "'6.04"
$ws.Range("E21").Value = "  -1.44%  "
$ws.Range("D22").Value = "'67.58"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").Value = "'236.67"
$ws.Range("E23").Value = "  -0.73%  "
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  +2.09%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'25.46"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").Value = "'167.27"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "'2.06"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").Value = "'9.09"
$ws.Range("E31").Value = "  -0.98%  "
$ws.Range("D32").Value = "'33.00"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("D34").Value = "'4.82"
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("D35").Value = "'5.00"
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("D36").Value = "'17.41"
$ws.Range("E36").Value = "  -5.15%  "
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "'0.0691"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").Value = "'1.75"
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("D42").Value = "'2.72"
$ws.Range("E42").Value = "  -1.30%  "
$ws.Range("D43").Value = "2.012.37"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("E44").Value = "  -2.19%  "
$ws.Range("D45").Value = "'18.07"
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("D46").Value = "'10.00"
$ws.Range("E46").Value = "  -2.21%  "
$ws.Range("D47").Value = "'2.11"
$ws.Range("E47").Value = "  -2.79%  "
$ws.Range("D48").Value = "'2.78"
$ws.Range("E48").Value = "  -2.00%  "
$ws.Range("D49").Value = "'2.93"
$ws.Range("E49").Value = "  +7.30%  "
$ws.Range("D50").Value = "'54.02"
$ws.Range("E50").Value = "  -1.02%  "
$ws.Range("D51").Value = "2.532.09"
$ws.Range("E51").Value = "  -0.24%  "
